# Apply the "add academic_id / sem_id / course_id columns" edit described
# in the commit diff to questions_checkbox_template.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Insert three new columns at the front of the sheet (A:C). Everything
#    that used to live in columns A:P shifts right to D:S, shared-string
#    references and formatting follow automatically.
$ws.Columns("A:C").Insert()

# 2. Populate the three new columns with their header + data values.
$ws.Range("A1").Value = "academic_id"
$ws.Range("B1").Value = "sem_id"
$ws.Range("C1").Value = "course_id"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 17

# 3. Match the authored column widths: A and B become fixed-width (13
#    characters). Using 12.2 here lands exactly on a raw width of 13 once
#    Excel's pixel-padding rounding is applied (matches the target XML).
$ws.Columns("A:B").ColumnWidth = 12.2

# 4. The defined names (_FilterDatabase / question_type / type) pointed at
#    column M, which is now column P after the 3-column insert — update
#    them explicitly since Insert() does not retarget defined names.
$ws.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$P`$3:`$P`$3"
$wb.Names.Item("question_type").RefersTo = "=Sheet1!`$P`$3:`$P`$3"
$wb.Names.Item("type").RefersTo = "=Sheet1!`$P`$3:`$P`$6"

# 5. Restore the sheet selection to match the authored workbook (whole
#    column C selected, active cell C1). Wrapped in [void] so the boolean
#    return value of Select() isn't echoed to the output stream.
[void]$ws.Range("C1:C1048576").Select()
